$d = $word.ActiveDocument

# --- 1. "ServiceNow Function V2.0.0" -> "ServiceNow App V2.0.0" ------------
# The title run " Function V" becomes " " + "App" + " V" (i.e. "Function" -> "App").
$d.Content.Find.Execute("Function", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "App", 2)

# --- 2. Turn on distinct odd/even-page footers -----------------------------
# Previously only "default" + "first" footers existed (footer1.xml / footer2.xml).
# Enabling OddAndEvenPagesHeaderFooter makes Word materialize a third, blank
# "even" footer and renumber the footer/header relationship ids, which also
# shuffles footer1.xml -> becomes the new blank even footer, the old
# footer1.xml ("Page X") content becomes footer2.xml (default), and the old
# footer2.xml (copyright) content becomes footer3.xml (first).
$sec = $d.Sections.Item(1)
$sec.PageSetup.OddAndEvenPagesHeaderFooter = $true

# Touch the (now existing) even-page footer so the engine actually materializes
# the new footer part/relationship/sectPr reference instead of just flipping
# the section flag.
$evenFooter = $sec.Footers.Item(3)
$evenFooter.Range.Text = ""
